# Weekly CompStat report refresh: new Police Commissioner, updated volume/week
# number, updated reporting dates, and a fresh week of crime-complaint figures
# for the 7th Precinct table (rows 15-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead updates -------------------------------------------------

# Police Commissioner name (M6, merged M6:N6)
$ws.Range("M6").Value = "Edward A. Caban"

# "Volume 30   Number  26" -> "...  27" (C8, merged C8:L8)
$ws.Range("C8").Value = "Volume 30   Number  27"

# "Report Covering the Week  6/26/2023  Through  7/2/2023" -> next week (C9, merged C9:L9)
$ws.Range("C9").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# --- Crime-complaints table (rows 15-30) --------------------------------
#
# Cells in this table are either plain numbers (counts / % changes) or, when
# a category had zero incidents in both periods, a placeholder pair of
# shared strings: "0" (count) / "***.*" (undefined % change). Several cells
# this week flip between those two representations, which also means their
# cell style flips between the "text" style (14) and the matching numeric
# style (15 count / 16 percent). A plain Value assignment can't change a
# cell from numeric style to text style (it keeps re-inferring a number), so
# Convert-ToTextCell stages the value under a Text number format and then
# copies formatting (PasteSpecial formats-only) from a donor cell that
# already carries the right target style, and Convert-ToNumCell does the
# mirror image for text -> number.

function Set-NumCell($ref, $val) {
    $ws.Range($ref).Value = $val
}

function Convert-ToNumCell($ref, $val, $donor) {
    $ws.Range($donor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $val
}

function Convert-ToTextCell($ref, $text, $donor) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range($donor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

function Set-NumCell($ref, $val) {
    $ws.Range($ref).Value = $val
}

function Convert-ToNumCell($ref, $val, $donor) {
    $ws.Range($donor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $val
}

function Convert-ToTextCell($ref, $text, $donor) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range($donor).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

Convert-ToNumCell "C15" 1 "F14"
Convert-ToNumCell "D15" 1 "F14"
Convert-ToNumCell "E15" 0 "K14"
Convert-ToNumCell "F15" 1 "F14"
Set-NumCell "G15" 2
Set-NumCell "H15" -50
Set-NumCell "I15" 6
Set-NumCell "J15" 10
Set-NumCell "K15" -40
Set-NumCell "L15" -14.285714285714
Set-NumCell "M15" 50
Set-NumCell "N15" -33.333333333333
Set-NumCell "C16" 1
Set-NumCell "D16" 3
Set-NumCell "E16" -66.666666666666
Set-NumCell "F16" 12
Set-NumCell "G16" 12
Set-NumCell "H16" 0
Set-NumCell "I16" 66
Set-NumCell "J16" 89
Set-NumCell "K16" -25.842696629213
Set-NumCell "L16" 43.478260869565
Set-NumCell "M16" -13.157894736842
Set-NumCell "N16" -83.582089552238
Set-NumCell "C17" 3
Set-NumCell "D17" 8
Set-NumCell "E17" -62.5
Set-NumCell "F17" 19
Set-NumCell "G17" 18
Set-NumCell "H17" 5.555555555555
Set-NumCell "I17" 108
Set-NumCell "J17" 114
Set-NumCell "K17" -5.263157894736
Set-NumCell "L17" 17.391304347826
Set-NumCell "M17" 35
Set-NumCell "N17" -6.896551724137
Set-NumCell "C18" 2
Set-NumCell "D18" 1
Set-NumCell "E18" 100
Set-NumCell "F18" 14
Set-NumCell "G18" 8
Set-NumCell "H18" 75
Set-NumCell "I18" 77
Set-NumCell "J18" 81
Set-NumCell "K18" -4.938271604938
Set-NumCell "L18" 40
Set-NumCell "M18" 97.435897435897
Set-NumCell "N18" -58.378378378378
Set-NumCell "D19" 11
Set-NumCell "E19" 9.090909090909
Set-NumCell "F19" 51
Set-NumCell "G19" 49
Set-NumCell "H19" 4.081632653061
Set-NumCell "I19" 305
Set-NumCell "J19" 398
Set-NumCell "K19" -23.366834170854
Set-NumCell "L19" 22
Set-NumCell "M19" 125.925925925926
Set-NumCell "N19" 30.341880341880
Convert-ToNumCell "C20" 2 "F14"
Set-NumCell "E20" 0
Set-NumCell "F20" 5
Set-NumCell "G20" 9
Set-NumCell "H20" -44.444444444444
Set-NumCell "I20" 29
Set-NumCell "J20" 34
Set-NumCell "K20" -14.705882352941
Set-NumCell "L20" -19.444444444444
Set-NumCell "M20" 3.571428571428
Set-NumCell "N20" -84.974093264248
Set-NumCell "C21" 21
Set-NumCell "D21" 26
Set-NumCell "E21" -19.230769230769
Set-NumCell "F21" 103
Set-NumCell "G21" 98
Set-NumCell "H21" 5.102040816326
Set-NumCell "I21" 593
Set-NumCell "J21" 728
Set-NumCell "K21" -18.543956043956
Set-NumCell "L21" 21.765913757700
Set-NumCell "M21" 63.360881542699
Set-NumCell "N21" -48.164335664335
Convert-ToTextCell "C22" "0" "C14"
Set-NumCell "F22" 2
Convert-ToTextCell "G22" "0" "C14"
Convert-ToTextCell "H22" "***.*" "C14"
Set-NumCell "L22" 37.5
Set-NumCell "M22" 57.142857142857
Set-NumCell "C23" 7
Set-NumCell "D23" 4
Set-NumCell "E23" 75
Set-NumCell "F23" 16
Set-NumCell "G23" 15
Set-NumCell "H23" 6.666666666666
Set-NumCell "I23" 80
Set-NumCell "J23" 97
Set-NumCell "K23" -17.525773195876
Set-NumCell "L23" -16.666666666666
Set-NumCell "M23" 1.265822784810
Set-NumCell "C24" 29
Set-NumCell "D24" 31
Set-NumCell "E24" -6.451612903225
Set-NumCell "F24" 100
Set-NumCell "G24" 200
Set-NumCell "H24" -50
Set-NumCell "I24" 641
Set-NumCell "J24" 1305
Set-NumCell "K24" -50.881226053639
Set-NumCell "L24" 2.724358974358
Set-NumCell "M24" 65.206185567010
Set-NumCell "C25" 4
Set-NumCell "D25" 10
Set-NumCell "E25" -60
Set-NumCell "F25" 38
Set-NumCell "G25" 28
Set-NumCell "H25" 35.714285714285
Set-NumCell "I25" 233
Set-NumCell "J25" 234
Set-NumCell "K25" -0.427350427350
Set-NumCell "L25" 64.084507042253
Set-NumCell "M25" 30.898876404494
Convert-ToNumCell "C26" 1 "F14"
Convert-ToNumCell "D26" 2 "F14"
Convert-ToNumCell "E26" -50 "K14"
Convert-ToNumCell "F26" 1 "F14"
Set-NumCell "G26" 3
Set-NumCell "H26" -66.666666666666
Set-NumCell "I26" 12
Set-NumCell "J26" 17
Set-NumCell "K26" -29.411764705882
Set-NumCell "L26" 20
Convert-ToNumCell "C27" 1 "F14"
Convert-ToTextCell "D27" "0" "C14"
Convert-ToTextCell "E27" "***.*" "C14"
Set-NumCell "G27" 6
Set-NumCell "H27" -66.666666666666
Set-NumCell "I27" 28
Set-NumCell "K27" 16.666666666666
Set-NumCell "L27" 75
Convert-ToTextCell "C28" "0" "C14"
Set-NumCell "E28" -100
Set-NumCell "G28" 2
Set-NumCell "H28" 100
Set-NumCell "J28" 8
Set-NumCell "K28" -25
Set-NumCell "M28" 50
Convert-ToTextCell "C29" "0" "C14"
Set-NumCell "E29" -100
Set-NumCell "G29" 2
Set-NumCell "H29" 50
Set-NumCell "J29" 6
Set-NumCell "K29" -16.666666666666
Set-NumCell "M29" 66.666666666666
Convert-ToTextCell "D30" "0" "C14"
Convert-ToTextCell "E30" "***.*" "C14"
Set-NumCell "G30" 2

